# Apply the latest 1h crypto snapshot: update Price (D) and Volume/1h-change (E)
# columns for the affected rows. Source data is plain text (inlineStr) -- for the
# Price column, some new values read as valid numbers (e.g. "64.00"), which Excel
# would silently coerce to a Number and normalize (dropping the trailing zero, or
# introducing floating point noise). Briefly mark the cell as Text, assign the
# literal string, then clear the format again so the cell keeps its original
# (unstyled) look -- only the text content changes, exactly like the source diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.191.15'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -0.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.216.27'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '301.22'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.97%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '88.43'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -5.47%  '
$ws.Range("E7").Value = '  -3.27%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  -6.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '33.36'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -3.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0779'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -3.30%  '
$ws.Range("E12").Value = '  -1.14%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '6.85'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -4.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.555.10'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.293.34'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.796'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.76%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.05'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.73%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.982.57'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0899'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -6.74%  '
$ws.Range("E20").Value = '  -5.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.13'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -7.28%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '64.00'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '231.84'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.13%  '
$ws.Range("E24").Value = '  -2.94%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.89'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -5.66%  '
$ws.Range("E27").Value = '  +1.58%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.35'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -4.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '35.98'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -9.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.31'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -3.52%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.56'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -5.71%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '145.74'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0751'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -5.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.91'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.09%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.106'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.00%  '
$ws.Range("E37").Value = '  -3.99%  '
$ws.Range("E38").Value = '  -2.94%  '
$ws.Range("E39").Value = '  -0.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.17'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -8.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.59'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.02%  '
$ws.Range("E42").Value = '  -4.16%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.730.19'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +1.75%  '
$ws.Range("E45").Value = '  +4.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '77.42'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -5.97%  '
$ws.Range("E47").Value = '  -6.38%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '94.19'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -5.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '66.16'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.59'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -6.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.436.25'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -1.02%  '
